# Add columns I (I0) and J (IF) to the sheet, following the style of the
# existing header cell H1, and fill in the per-row data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy H1's formatting (bold font + border) onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Data rows (2-51) ---
$I = @(8,9,9,7,8,7,7,7,9,8,8,9,7,8,7,8,7,9,9,9,7,7,8,6,7,9,11,7,8,7,10,7,7,7,8,9,8,7,6,7,6,7,6,7,7,5,5,6,5,3)
$J = @(8,9,9,8,8,7,7,7,9,8,8,9,7,8,7,8,7,9,9,9,7,7,8,6,7,9,11,7,8,7,10,7,7,8,8,9,8,7,6,8,6,7,6,7,7,5,5,6,5,3)

for ($i = 0; $i -lt $I.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $I[$i]
    $ws.Cells.Item($row, 10).Value = $J[$i]
}
